$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")

# Read current text (use Value2, which reads back the real string in this
# runtime) and patch only the two changed lines, leaving everything else
# - including blank lines, emoji and accents - untouched.
$currentText = $cellA1.Value2
$updatedText = $currentText -replace [regex]::Escape("1000 Bs = 11.24 = 45606.74 pesos"), "1000 Bs = 11.2 = 45413.66 pesos"
$updatedText = $updatedText -replace [regex]::Escape("45606.74 pesos = 11.21 = 982.97 Bs"), "45413.66 pesos = 11.14 = 956.76 Bs"

$cellA1.Value = $updatedText

# --- Sheet "tasas": update the rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 89.3
$wsTasas.Range("O10").Value = 4055.44
$wsTasas.Range("N12").Value = 4078.3
$wsTasas.Range("O12").Value = 85.92
